# Natmi following Dr Hou advice
# Update rows 2-3 and add rows 4-5 with recomputed ligand/receptor
# statistics after merging FAPs + sCs clusters (per Dr Hou advice).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "FAPs"
$ws.Cells.Item(2, 2).Value = "Ngf"
$ws.Cells.Item(2, 3).Value = "Ntrk1"
$ws.Cells.Item(2, 4).Value = "FAPs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.7246536666666668
$ws.Cells.Item(2, 8).Value = 2.173961
$ws.Cells.Item(2, 9).Value = 0.1791272621505297
$ws.Cells.Item(2, 10).Value = 0.1791272621505298
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.08849299999999999
$ws.Cells.Item(2, 14).Value = 0.265479
$ws.Cells.Item(2, 15).Value = 0.9052313210944106
$ws.Cells.Item(2, 16).Value = 0.9052313210944106
$ws.Cells.Item(2, 17).Value = 0.06412677692433333
$ws.Cells.Item(2, 18).Value = 0.577140992319
$ws.Cells.Item(2, 19).Value = 0.1621516081605489
$ws.Cells.Item(2, 20).Value = 0.1621516081605489

# Row 3
$ws.Cells.Item(3, 1).Value = "FAPs"
$ws.Cells.Item(3, 2).Value = "Ngf"
$ws.Cells.Item(3, 3).Value = "Ntrk1"
$ws.Cells.Item(3, 4).Value = "sCs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.7246536666666668
$ws.Cells.Item(3, 8).Value = 2.173961
$ws.Cells.Item(3, 9).Value = 0.1791272621505297
$ws.Cells.Item(3, 10).Value = 0.1791272621505298
$ws.Cells.Item(3, 11).Value = 1
$ws.Cells.Item(3, 12).Value = 0.3333333333333333
$ws.Cells.Item(3, 13).Value = 0.009264333333333334
$ws.Cells.Item(3, 14).Value = 0.027793
$ws.Cells.Item(3, 15).Value = 0.09476867890558938
$ws.Cells.Item(3, 16).Value = 0.09476867890558936
$ws.Cells.Item(3, 17).Value = 0.006713433119222224
$ws.Cells.Item(3, 18).Value = 0.06042089807300001
$ws.Cells.Item(3, 19).Value = 0.01697565398998089
$ws.Cells.Item(3, 20).Value = 0.01697565398998089

# Row 4
$ws.Cells.Item(4, 1).Value = "sCs"
$ws.Cells.Item(4, 2).Value = "Ngf"
$ws.Cells.Item(4, 3).Value = "Ntrk1"
$ws.Cells.Item(4, 4).Value = "FAPs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 3.320814666666667
$ws.Cells.Item(4, 8).Value = 9.962444
$ws.Cells.Item(4, 9).Value = 0.8208727378494701
$ws.Cells.Item(4, 10).Value = 0.8208727378494702
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.08849299999999999
$ws.Cells.Item(4, 14).Value = 0.265479
$ws.Cells.Item(4, 15).Value = 0.9052313210944106
$ws.Cells.Item(4, 16).Value = 0.9052313210944106
$ws.Cells.Item(4, 17).Value = 0.2938688522973333
$ws.Cells.Item(4, 18).Value = 2.644819670676
$ws.Cells.Item(4, 19).Value = 0.7430797129338617
$ws.Cells.Item(4, 20).Value = 0.7430797129338618

# Row 5
$ws.Cells.Item(5, 1).Value = "sCs"
$ws.Cells.Item(5, 2).Value = "Ngf"
$ws.Cells.Item(5, 3).Value = "Ntrk1"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 3.320814666666667
$ws.Cells.Item(5, 8).Value = 9.962444
$ws.Cells.Item(5, 9).Value = 0.8208727378494701
$ws.Cells.Item(5, 10).Value = 0.8208727378494702
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.009264333333333334
$ws.Cells.Item(5, 14).Value = 0.027793
$ws.Cells.Item(5, 15).Value = 0.09476867890558938
$ws.Cells.Item(5, 16).Value = 0.09476867890558936
$ws.Cells.Item(5, 17).Value = 0.03076513401022223
$ws.Cells.Item(5, 18).Value = 0.276886206092
$ws.Cells.Item(5, 19).Value = 0.07779302491560848
$ws.Cells.Item(5, 20).Value = 0.07779302491560848

Write-Output "Updated rows 2-5 of Sheet1"
